$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000496.2
$ws.Range("J2").Value = 762.5
$ws.Range("L2").Value = 762.5
$ws.Range("N2").Value = -988.5
$ws.Range("H9").Value = 1429764.9
$ws.Range("I9").Value = 2501300.5
$ws.Range("J9").Value = 1050.6666
$ws.Range("K9").Value = 2501300.5
$ws.Range("L9").Value = 1050.6666
$ws.Range("M9").Value = -2501131.5
$ws.Range("N9").Value = -1388.6666
$ws.Range("H54").Value = 14530.4
$ws.Range("I54").Value = 14530.4
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 14530.4
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -14044.4
$ws.Range("N54").Value = $null
$ws.Range("H116").Value = 16977.6
$ws.Range("J116").Value = 14963.667
$ws.Range("L116").Value = 14963.667
$ws.Range("N116").Value = -21847.667
$ws.Range("H127").Value = 1325.8
$ws.Range("I127").Value = 1462
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 4386
$ws.Range("L127").Value = 300
$ws.Range("M127").Value = 574
$ws.Range("N127").Value = -10220
$ws.Range("H129").Value = 749.9231
$ws.Range("I129").Value = 515
$ws.Range("J129").Value = 1533
$ws.Range("K129").Value = 1545
$ws.Range("L129").Value = 4599
$ws.Range("M129").Value = 3455
$ws.Range("N129").Value = -14599
$ws.Range("H138").Value = 3658.2744
$ws.Range("I138").Value = 3354.889
$ws.Range("J138").Value = 3723.2856
$ws.Range("K138").Value = 10064.667
$ws.Range("L138").Value = 11169.8568
$ws.Range("M138").Value = -4924.667000000001
$ws.Range("N138").Value = -21449.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2174.3076
$ws.Range("I45").Value = 1351.8889
$ws.Range("J45").Value = 4024.75
$ws.Range("K45").Value = 1351.8889
$ws.Range("L45").Value = 4024.75
$ws.Range("M45").Value = -974.8888999999999
$ws.Range("N45").Value = -4778.75
$ws.Range("H74").Value = 5220.6523
$ws.Range("I74").Value = 4470.143
$ws.Range("J74").Value = 6388.1113
$ws.Range("K74").Value = 4470.143
$ws.Range("L74").Value = 6388.1113
$ws.Range("M74").Value = -3596.143
$ws.Range("N74").Value = -8136.1113
$ws.Range("H77").Value = 5220.6523
$ws.Range("I77").Value = 4470.143
$ws.Range("J77").Value = 6388.1113
$ws.Range("K77").Value = 22350.715
$ws.Range("L77").Value = 31940.5565
$ws.Range("M77").Value = -17982.715
$ws.Range("N77").Value = -40676.5565
$ws.Range("H102").Value = 4000.2778
$ws.Range("I102").Value = 2800.3333
$ws.Range("K102").Value = 2800.3333
$ws.Range("M102").Value = -1178.3333
$ws.Range("H122").Value = 3585759.2
$ws.Range("I122").Value = 3969754.8
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 11909264.4
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -11906814.4
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1215.909
$ws.Range("I20").Value = 838.1429000000001
$ws.Range("J20").Value = 1877
$ws.Range("K20").Value = 838.1429000000001
$ws.Range("L20").Value = 1877
$ws.Range("M20").Value = -591.1429000000001
$ws.Range("N20").Value = -2371
$ws.Range("H70").Value = 299988
$ws.Range("J70").Value = 299988
$ws.Range("L70").Value = 299988
$ws.Range("N70").Value = -300574
$ws.Range("H73").Value = 299988
$ws.Range("J73").Value = 299988
$ws.Range("L73").Value = 299988
$ws.Range("N73").Value = -302016
$ws.Range("H87").Value = 81400
$ws.Range("J87").Value = 81400
$ws.Range("L87").Value = 81400
$ws.Range("N87").Value = -83896
$ws.Range("H90").Value = 81400
$ws.Range("J90").Value = 81400
$ws.Range("L90").Value = 244200
$ws.Range("N90").Value = -256680
$ws.Range("H105").Value = 2932.889
$ws.Range("I105").Value = 2166
$ws.Range("J105").Value = 4466.6665
$ws.Range("K105").Value = 2166
$ws.Range("L105").Value = 4466.6665
$ws.Range("M105").Value = -419
$ws.Range("N105").Value = -7960.6665
$ws.Range("H107").Value = 9656.714
$ws.Range("J107").Value = 9832.933999999999
$ws.Range("L107").Value = 9832.933999999999
$ws.Range("N107").Value = -13672.934

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2330.889
$ws.Range("I31").Value = 1130.1305
$ws.Range("J31").Value = 4455.3076
$ws.Range("K31").Value = 1130.1305
$ws.Range("L31").Value = 4455.3076
$ws.Range("M31").Value = -835.1305
$ws.Range("N31").Value = -5045.3076
$ws.Range("H34").Value = 2330.889
$ws.Range("I34").Value = 1130.1305
$ws.Range("J34").Value = 4455.3076
$ws.Range("K34").Value = 1130.1305
$ws.Range("L34").Value = 4455.3076
$ws.Range("M34").Value = -928.1305
$ws.Range("N34").Value = -4859.3076
$ws.Range("H93").Value = 11420.667
$ws.Range("I93").Value = 12223.25
$ws.Range("K93").Value = 12223.25
$ws.Range("M93").Value = -10351.25
$ws.Range("H103").Value = 24999.666
$ws.Range("I103").Value = 22499.5
$ws.Range("K103").Value = 22499.5
$ws.Range("M103").Value = -21327.5
$ws.Range("H134").Value = 2187.875
$ws.Range("I134").Value = 2300.5833
$ws.Range("J134").Value = 1849.75
$ws.Range("K134").Value = 6901.749899999999
$ws.Range("L134").Value = 5549.25
$ws.Range("M134").Value = -4366.749899999999
$ws.Range("N134").Value = -10619.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 182.375
$ws.Range("I33").Value = 201.28572
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 1207.71432
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -924.71432
$ws.Range("N33").Value = -866
$ws.Range("H46").Value = 22424608
$ws.Range("I46").Value = 23232812
$ws.Range("K46").Value = 69698436
$ws.Range("M46").Value = -69698345
$ws.Range("H121").Value = 823.0769
$ws.Range("I121").Value = 753.125
$ws.Range("J121").Value = 935
$ws.Range("K121").Value = 2259.375
$ws.Range("L121").Value = 2805
$ws.Range("M121").Value = -949.375
$ws.Range("N121").Value = -5425
$ws.Range("H131").Value = 2987762.5
$ws.Range("I131").Value = 333708.34
$ws.Range("J131").Value = 3924487.5
$ws.Range("K131").Value = 1001125.02
$ws.Range("L131").Value = 11773462.5
$ws.Range("M131").Value = -996085.02
$ws.Range("N131").Value = -11783542.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 6000
$ws.Range("J49").Value = 6000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6368
$ws.Range("H113").Value = 28577646
$ws.Range("I113").Value = 71430440
$ws.Range("J113").Value = 9118.477000000001
$ws.Range("K113").Value = 71430440
$ws.Range("L113").Value = 9118.477000000001
$ws.Range("M113").Value = -71428270
$ws.Range("N113").Value = -13458.477

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7783.115
$ws.Range("I40").Value = 7128.8125
$ws.Range("J40").Value = 8830
$ws.Range("K40").Value = 7128.8125
$ws.Range("L40").Value = 8830
$ws.Range("M40").Value = -6992.8125
$ws.Range("N40").Value = -9102
$ws.Range("H46").Value = 3974.963
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3974.963
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3974.963
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -4350.963
$ws.Range("H88").Value = 69449
$ws.Range("J88").Value = 69449
$ws.Range("L88").Value = 69449
$ws.Range("N88").Value = -70305
$ws.Range("H91").Value = 69449
$ws.Range("J91").Value = 69449
$ws.Range("L91").Value = 69449
$ws.Range("N91").Value = -72413
$ws.Range("H100").Value = 4300.2
$ws.Range("I100").Value = 2389.2222
$ws.Range("J100").Value = 7166.6665
$ws.Range("K100").Value = 2389.2222
$ws.Range("L100").Value = 7166.6665
$ws.Range("M100").Value = -1848.2222
$ws.Range("N100").Value = -8248.666499999999
$ws.Range("H122").Value = 4674.722
$ws.Range("I122").Value = 4489.8
$ws.Range("K122").Value = 13469.4
$ws.Range("M122").Value = -11019.4
$ws.Range("H132").Value = 5510.591
$ws.Range("I132").Value = 5249
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 15747
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -13217
$ws.Range("N132").Value = -24260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4375.25
$ws.Range("I4").Value = 2901
$ws.Range("J4").Value = 5849.5
$ws.Range("K4").Value = 2901
$ws.Range("L4").Value = 5849.5
$ws.Range("M4").Value = -2788
$ws.Range("N4").Value = -6075.5
$ws.Range("H70").Value = 41192.223
$ws.Range("J70").Value = 41192.223
$ws.Range("L70").Value = 41192.223
$ws.Range("N70").Value = -41822.223
$ws.Range("H73").Value = 41192.223
$ws.Range("J73").Value = 41192.223
$ws.Range("L73").Value = 41192.223
$ws.Range("N73").Value = -43376.223
$ws.Range("H107").Value = 2425.6
$ws.Range("I107").Value = 2907
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 8721
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -6801
$ws.Range("N107").Value = -5340
$ws.Range("H122").Value = 7824.875
$ws.Range("I122").Value = 6520
$ws.Range("J122").Value = 8418
$ws.Range("K122").Value = 19560
$ws.Range("L122").Value = 25254
$ws.Range("M122").Value = -17110
$ws.Range("N122").Value = -30154
$ws.Range("H132").Value = 1828.8214
$ws.Range("I132").Value = 1638.5652
$ws.Range("J132").Value = 2704
$ws.Range("K132").Value = 4915.6956
$ws.Range("L132").Value = 8112
$ws.Range("M132").Value = -2385.6956
$ws.Range("N132").Value = -13172
